# The "Classroom" sheet's AutoTest Setup block contained two now-obsolete
# rows describing the (LinkedList) git-clone / setup step and the
# "Main Output" step. Remove them the same way a user would in Excel:
# select rows 3:4 and delete the entire rows, which shifts everything
# below up by two rows and recalculates the trailing SUM().

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classroom")

$ws.Rows("3:4").Delete()

# Leave the selection on the row that slid into the old row 3's place,
# matching how Excel leaves the selection after a row deletion.
$ws.Range("A3:XFD3").Select()
